$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.456.22'
$ws.Range("E2").Value = '  +1.54%  '
$ws.Range("D3").Value = '1.806.66'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.81'
$ws.Range("E5").Value = '  -2.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9989'
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4463'
$ws.Range("E7").Value = '  +5.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3780'
$ws.Range("E8").Value = '  +7.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.66'
$ws.Range("E9").Value = '  -2.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.151'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07517'
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.62'
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9998'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.634'
$ws.Range("E14").Value = '  +4.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.307'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = '1.802.88'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001093'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06805'
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '80.71'
$ws.Range("E19").Value = '  -1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9993'
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.65'
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.339'
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("D23").Value = '28.437.87'
$ws.Range("E23").Value = '  +1.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.83'
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.406'
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.370'
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.54'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.05'
$ws.Range("E28").Value = '  -1.41%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.358'
$ws.Range("E29").Value = '  -4.64%  '
$ws.Range("B30").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C30").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D30").Value = '2.006.77'
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '132.68'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.260'
$ws.Range("E32").Value = '  -3.59%  '
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.004'
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.828'
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.09343'
$ws.Range("E35").Value = '  +1.87%  '
$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2277'
$ws.Range("E36").Value = '  +4.73%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '12.18'
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06358'
$ws.Range("E38").Value = '  +1.01%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02345'
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6598'
$ws.Range("E40").Value = '  -1.45%  '
$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.171'
$ws.Range("E41").Value = '  -1.18%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.211'
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("B43").Value = 'WEMIXTOKEN'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.452'
$ws.Range("E43").Value = '  -3.87%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.105'
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9986'
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.81'
$ws.Range("E46").Value = '  -3.24%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6088'
$ws.Range("E47").Value = '  -0.96%  '
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.814'
$ws.Range("E48").Value = '  -1.60%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '128.51'
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.036'
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07089'
$ws.Range("E51").Value = '  -0.40%  '
